# Fruta / hortaliza, semanal
#
# Inserts two new weekly price rows for "Comercializadora del Agro de
# Limarí" / "Ají" (Americana (o), Primera & Segunda) dated 2021-11-04
# (serial 44504) right after the existing row for 2021-09-30 (row 144),
# shifting every subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 145-146; everything that was on/after row 145
# (old rows 145..176) shifts down to 147..178.
$ws.Rows("145:146").Insert()

# --- New row 145: Americana (o) / Primera, 4-nov-2021 ---
$ws.Range("A145").Value = 2
$ws.Range("B145").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44504
$ws.Range("E145").Value = 4
$ws.Range("F145").Value = 100112021
$ws.Range("G145").Value = "Ají"
$ws.Range("H145").Value = "Americana (o)"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 240
$ws.Range("K145").Value = 30000
$ws.Range("L145").Value = 35000
$ws.Range("M145").Value = 32500
$ws.Range("N145").Value = "$/caja 25 kilos"
$ws.Range("O145").Value = "Provincia de Limarí"
$ws.Range("P145").Value = 1300
$ws.Range("Q145").Value = 25
$ws.Range("R145").Value = "Hortaliza"

# --- New row 146: Americana (o) / Segunda, 4-nov-2021 ---
$ws.Range("A146").Value = 2
$ws.Range("B146").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44504
$ws.Range("E146").Value = 4
$ws.Range("F146").Value = 100112021
$ws.Range("G146").Value = "Ají"
$ws.Range("H146").Value = "Americana (o)"
$ws.Range("I146").Value = "Segunda"
$ws.Range("J146").Value = 100
$ws.Range("K146").Value = 20000
$ws.Range("L146").Value = 25000
$ws.Range("M146").Value = 22500
$ws.Range("N146").Value = "$/caja 25 kilos"
$ws.Range("O146").Value = "Provincia de Limarí"
$ws.Range("P146").Value = 900
$ws.Range("Q146").Value = 25
$ws.Range("R146").Value = "Hortaliza"
